# Update workbook/sheet title to reflect the new "through" date
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2021-12-04"

# Row 14 - December (through 12-0x) month row: update label and figures,
# and add new no_arrest_made / arrest_rate values for the 2015 (B/C/D) group
$ws.Range("A14").Value = "December (through 12-04)"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 0.5
$ws.Range("D14").NumberFormat = "0.0%"

$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 12
$ws.Range("G14").Value = 0.0769

$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 15
$ws.Range("J14").Value = 0.0625

$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 10
$ws.Range("M14").Value = 0.0909

$ws.Range("O14").Value = 4
$ws.Range("R14").Value = 19
$ws.Range("U14").Value = 31

# Row 15 - Total row: update totals/rates to include new December data
$ws.Range("C15").Value = 259
$ws.Range("D15").Value = 0.116

$ws.Range("F15").Value = 515
$ws.Range("G15").Value = 0.1059

$ws.Range("I15").Value = 773
$ws.Range("J15").Value = 0.0765

$ws.Range("L15").Value = 618
$ws.Range("M15").Value = 0.1082

$ws.Range("O15").Value = 484
$ws.Range("P15").Value = 0.1004

$ws.Range("R15").Value = 1219
$ws.Range("S15").Value = 0.0499

$ws.Range("U15").Value = 1576
$ws.Range("V15").Value = 0.0591
